$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing block
# (old rows 236-250 shift down to 237-251); fill the freed row 236 with
# the new record's data.
$ws.Rows.Item(236).Insert()

$ws.Cells.Item(236, 1).Value = 5
$ws.Cells.Item(236, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(236, 3).Value = "Maule"
$ws.Cells.Item(236, 4).Value = 44706
$ws.Cells.Item(236, 5).Value = 7
$ws.Cells.Item(236, 6).Value = 100112009
$ws.Cells.Item(236, 7).Value = "Acelga"
$ws.Cells.Item(236, 8).Value = "Sin especificar"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 500
$ws.Cells.Item(236, 11).Value = 3000
$ws.Cells.Item(236, 12).Value = 3000
$ws.Cells.Item(236, 13).Value = 3000
$ws.Cells.Item(236, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(236, 15).Value = "Región del Maule"
$ws.Cells.Item(236, 16).Value = 750
$ws.Cells.Item(236, 17).Value = 4
$ws.Cells.Item(236, 18).Value = "Hortaliza"
